# "Add files via upload" - refresh of the dispatch-sheet creation dates.
#
# What changes in the workbook:
#   1. The "创建日期" (creation date) column P, rows 2-52, gets new date
#      serials written to it.
#   2. Before the edit, P2 carried a one-off cell format (no fill) while
#      P3:P52 shared a second, near-identical date format (white fill).
#      The refreshed file collapses these into a single shared format, so
#      we copy P3's format onto P2 first - after that every cell in the
#      column is formatted identically, matching the de-duplicated style
#      table in the target file.
#   3. The active selection moves from L54 to P54.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unify P2's format with the rest of the date column -----------------
$ws.Range("P3").Copy()
$ws.Range("P2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Write the refreshed date serials for P2:P52 -------------------------
$newDates = @(
    45809, 45810, 45819, 45819, 45819, 45819, 45819, 45819, 45819, 45819,
    45819, 45819, 45819, 45819, 45819, 45819, 45819, 45819, 45819, 45819,
    45819, 45819, 45819, 45819, 45819, 45820, 45820, 45820, 45820, 45820,
    45820, 45820, 45820, 45820, 45820, 45820, 45820, 45820, 45820, 45820,
    45820, 45820, 45820, 45820, 45820, 45820, 45820, 45820, 45820, 45820,
    45820
)

$startRow = 2
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 16).Value = $newDates[$i]
}

# --- Move the selection from L54 to P54 -----------------------------------
$ws.Range("P54").Select()
